$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data (OPPORTUNITY ID column keeps text storage, like the source data)
$ws.Range("A2").Value = "'1331979"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1331979"
$ws.Range("C2").Value = "Marketing Analyst Oncology and Specialities"
$ws.Range("D2").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "GSK Panama"

# Update row 3 data
$ws.Range("A3").Value = "'1331881"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331881"
$ws.Range("C3").Value = "Full stack mobile developer"
$ws.Range("D3").Value = "Kafr El-Shaikh, Qism Kafr El-Shaikh, Kafr el-Sheikh, Gharbia Governorate, Egypt"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Legend (حضانه ليجند)"

# Update row 4 data
$ws.Range("A4").Value = "'1326205"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326205"
$ws.Range("C4").Value = "Web & System Developer"
$ws.Range("D4").Value = "Nugegoda, Sri Lanka"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "57 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "All In One Holdings (PVT) Ltd"

# Remove rows 5, 6, 7 entirely (shrinks dimension to A1:H4)
$ws.Range("A5:H7").Delete()

# Adjust column widths to match new layout.
# NOTE: this COM bridge's ColumnWidth setter round-trips through a fixed
# +0.8333333333333334 (5/6 character) padding, so the raw OOXML <col width>
# ends up 0.8333... higher than the value assigned. Pre-subtract that
# constant so the persisted width lands exactly on the intended value.
$colPad = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 46 - $colPad
$ws.Columns.Item(4).ColumnWidth = 82 - $colPad
$ws.Columns.Item(6).ColumnWidth = 16 - $colPad
$ws.Columns.Item(7).ColumnWidth = 15 - $colPad
$ws.Columns.Item(8).ColumnWidth = 32 - $colPad
